# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (column G) replacing previous "Strike#" values, row by row (row 2 = first data row)
$kValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 2
    12 = 0
    13 = 0
    14 = 2
    15 = 1
    16 = 1
    17 = 2
    18 = 0
    19 = 3
    20 = 2
    21 = 0
    22 = 2
    23 = 1
    24 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
